# Update crypto price/volume table with freshly scraped figures.
# Values that look like plain decimal numbers (e.g. "1.00", "0.999") are
# prefixed with a literal leading apostrophe so Excel stores them as text
# (preserving exact formatting such as trailing zeros) instead of
# re-interpreting them as floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.348.41'
$ws.Range("E2").Value = '  +4.11%  '
$ws.Range("D3").Value = '3.484.91'
$ws.Range("E3").Value = '  +3.49%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''584.54'
$ws.Range("E5").Value = '  +2.29%  '
$ws.Range("D6").Value = '''147.39'
$ws.Range("E6").Value = '  +6.73%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +1.24%  '
$ws.Range("D9").Value = '''7.68'
$ws.Range("E9").Value = '  -0.40%  '
$ws.Range("D10").Value = '''0.126'
$ws.Range("E10").Value = '  +4.13%  '
$ws.Range("D11").Value = '''0.396'
$ws.Range("E11").Value = '  +3.70%  '
$ws.Range("D12").Value = '4.081.50'
$ws.Range("E12").Value = '  +3.58%  '
$ws.Range("D13").Value = '''29.64'
$ws.Range("E13").Value = '  +5.34%  '
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").Value = '3.478.87'
$ws.Range("E15").Value = '  +3.32%  '
$ws.Range("D16").Value = '''0.0000173'
$ws.Range("E16").Value = '  +3.23%  '
$ws.Range("D17").Value = '63.308.48'
$ws.Range("E17").Value = '  +3.87%  '
$ws.Range("D18").Value = '''6.29'
$ws.Range("E18").Value = '  +3.10%  '
$ws.Range("D19").Value = '''14.38'
$ws.Range("E19").Value = '  +6.48%  '
$ws.Range("D20").Value = '''9.36'
$ws.Range("E20").Value = '  +4.70%  '
$ws.Range("D21").Value = '''391.25'
$ws.Range("E21").Value = '  +1.98%  '
$ws.Range("D22").Value = '''0.563'
$ws.Range("E22").Value = '  +2.42%  '
$ws.Range("D23").Value = '''75.30'
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").Value = '''0.0000118'
$ws.Range("E25").Value = '  +7.60%  '
$ws.Range("D26").Value = '3.628.88'
$ws.Range("E26").Value = '  +3.59%  '
$ws.Range("E27").Value = '  -2.46%  '
$ws.Range("D28").Value = '''7.83'
$ws.Range("E28").Value = '  +10.10%  '
$ws.Range("D29").Value = '''1.00'
$ws.Range("D30").Value = '''8.26'
$ws.Range("E30").Value = '  +4.38%  '
$ws.Range("E31").Value = '  +1.89%  '
$ws.Range("D32").Value = '''1.43'
$ws.Range("E32").Value = '  +7.29%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").Value = '''23.79'
$ws.Range("E34").Value = '  +3.53%  '
$ws.Range("D35").Value = '''32.48'
$ws.Range("E35").Value = '  +26.66%  '
$ws.Range("D36").Value = '''5.32'
$ws.Range("E36").Value = '  +8.34%  '
$ws.Range("D37").Value = '''7.12'
$ws.Range("E37").Value = '  +4.49%  '
$ws.Range("D38").Value = '''171.84'
$ws.Range("E38").Value = '  +2.78%  '
$ws.Range("E39").Value = '  +9.16%  '
$ws.Range("D40").Value = '3.521.30'
$ws.Range("E40").Value = '  +3.45%  '
$ws.Range("E41").Value = '  +1.53%  '
$ws.Range("E42").Value = '  +4.60%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '''4.50'
$ws.Range("E43").Value = '  +3.74%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = '''42.47'
$ws.Range("E44").Value = '  +0.54%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").Value = '''1.73'
$ws.Range("E45").Value = '  +6.60%  '
$ws.Range("D46").Value = '''1.21'
$ws.Range("E46").Value = '  +9.57%  '
$ws.Range("D47").Value = '2.623.35'
$ws.Range("E47").Value = '  +7.61%  '
$ws.Range("D48").Value = '''23.69'
$ws.Range("E48").Value = '  +7.45%  '
$ws.Range("D49").Value = '''2.29'
$ws.Range("E49").Value = '  +16.57%  '
$ws.Range("D50").Value = '''6.74'
$ws.Range("E50").Value = '  +1.84%  '
$ws.Range("E51").Value = '  +5.25%  '
